$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2352941176470588
$ws.Range("C2").Value = 0.5032679738562091
$ws.Range("J2").Value = 0.01633986928104575
$ws.Range("P2").Value = 0.1568627450980392
$ws.Range("S2").Value = 0.08823529411764706
$ws.Range("B3").Value = 0.01273885350318471
$ws.Range("C3").Value = 0.006369426751592357
$ws.Range("J3").Value = 0.01273885350318471
$ws.Range("P3").Value = 0.7197452229299363
$ws.Range("S3").Value = 0.2484076433121019
$ws.Range("J4").Value = 0.0625
$ws.Range("P4").Value = 0.7083333333333334
$ws.Range("S4").Value = 0.2291666666666667
$ws.Range("B6").Value = 0.06707317073170732
$ws.Range("D6").Value = 0.01829268292682927
$ws.Range("F6").Value = 0.07926829268292683
$ws.Range("J6").Value = 0.2682926829268293
$ws.Range("O6").Value = 0.01219512195121951
$ws.Range("Q6").Value = 0.1890243902439024
$ws.Range("R6").Value = 0.06097560975609756
$ws.Range("S6").Value = 0.3048780487804878
$ws.Range("B7").Value = 0.1160714285714286
$ws.Range("F7").Value = 0.08035714285714286
$ws.Range("J7").Value = 0.1875
$ws.Range("Q7").Value = 0.1071428571428571
$ws.Range("R7").Value = 0.0625
$ws.Range("S7").Value = 0.4464285714285715
$ws.Range("B8").Value = 0.1349206349206349
$ws.Range("D8").Value = 0.03174603174603174
$ws.Range("E8").Value = 0.005291005291005291
$ws.Range("F8").Value = 0.04497354497354497
$ws.Range("J8").Value = 0.1428571428571428
$ws.Range("O8").Value = 0.02380952380952381
$ws.Range("Q8").Value = 0.1481481481481481
$ws.Range("R8").Value = 0.07407407407407407
$ws.Range("S8").Value = 0.3941798941798942
$ws.Range("B9").Value = 0.06870229007633588
$ws.Range("D9").Value = 0.03053435114503817
$ws.Range("F9").Value = 0.02290076335877863
$ws.Range("J9").Value = 0.1603053435114504
$ws.Range("Q9").Value = 0.2595419847328244
$ws.Range("R9").Value = 0.09923664122137404
$ws.Range("S9").Value = 0.3587786259541985
$ws.Range("B10").Value = 0.1331496786042241
$ws.Range("D10").Value = 0.02479338842975207
$ws.Range("E10").Value = 0.002754820936639119
$ws.Range("F10").Value = 0.05968778696051423
$ws.Range("J10").Value = 0.1092745638200184
$ws.Range("O10").Value = 0.01928374655647383
$ws.Range("Q10").Value = 0.224058769513315
$ws.Range("R10").Value = 0.1019283746556474
$ws.Range("S10").Value = 0.325068870523416
$ws.Range("G11").Value = 0.15625
$ws.Range("J11").Value = 0.09375
$ws.Range("K11").Value = 0.2135416666666667
$ws.Range("L11").Value = 0.5208333333333334
$ws.Range("S11").Value = 0.015625
$ws.Range("G12").Value = 0.6796116504854369
$ws.Range("J12").Value = 0.2330097087378641
$ws.Range("L12").Value = 0.01941747572815534
$ws.Range("S12").Value = 0.06796116504854369
$ws.Range("G13").Value = 0.72
$ws.Range("J13").Value = 0.2
$ws.Range("S13").Value = 0.08
$ws.Range("F15").Value = 0.01785714285714286
$ws.Range("H15").Value = 0.130952380952381
$ws.Range("I15").Value = 0.04166666666666666
$ws.Range("J15").Value = 0.4464285714285715
$ws.Range("K15").Value = 0.06547619047619048
$ws.Range("M15").Value = 0.01785714285714286
$ws.Range("O15").Value = 0.02380952380952381
$ws.Range("S15").Value = 0.2559523809523809
$ws.Range("F16").Value = 0.02127659574468085
$ws.Range("H16").Value = 0.175531914893617
$ws.Range("I16").Value = 0.04787234042553191
$ws.Range("J16").Value = 0.425531914893617
$ws.Range("K16").Value = 0.06382978723404255
$ws.Range("M16").Value = 0.01595744680851064
$ws.Range("N16").Value = 0.005319148936170213
$ws.Range("O16").Value = 0.09042553191489362
$ws.Range("S16").Value = 0.1542553191489362
$ws.Range("F17").Value = 0.01336898395721925
$ws.Range("H17").Value = 0.2058823529411765
$ws.Range("I17").Value = 0.09090909090909091
$ws.Range("J17").Value = 0.4411764705882353
$ws.Range("K17").Value = 0.06684491978609626
$ws.Range("M17").Value = 0.0160427807486631
$ws.Range("N17").Value = 0.00267379679144385
$ws.Range("O17").Value = 0.06951871657754011
$ws.Range("S17").Value = 0.09358288770053476
$ws.Range("F18").Value = 0.02380952380952381
$ws.Range("H18").Value = 0.1607142857142857
$ws.Range("I18").Value = 0.08333333333333333
$ws.Range("J18").Value = 0.4761904761904762
$ws.Range("K18").Value = 0.07142857142857142
$ws.Range("M18").Value = 0.01785714285714286
$ws.Range("O18").Value = 0.05357142857142857
$ws.Range("S18").Value = 0.1130952380952381
$ws.Range("F19").Value = 0.01900739176346357
$ws.Range("H19").Value = 0.2397043294614572
$ws.Range("I19").Value = 0.07180570221752904
$ws.Range("J19").Value = 0.4002111932418163
$ws.Range("K19").Value = 0.08764519535374868
$ws.Range("M19").Value = 0.01372756071805702
$ws.Range("N19").Value = 0.002111932418162619
$ws.Range("O19").Value = 0.06546990496304118
$ws.Range("S19").Value = 0.1003167898627244
